$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.051369071006775
$ws.Range("B1").Value = 1.104343891143799
$ws.Range("C1").Value = 0.8595868349075317
$ws.Range("D1").Value = 4.945044040679932
$ws.Range("E1").Value = 2.073944807052612
